# Auto-generated edit script: updates cached market-price / profit values
# in the Leve-profit tracker sheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR) to
# match a refreshed data pull. Values are plain numeric cell writes; no
# formulas are involved anywhere in this workbook.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1139.7
$ws.Range("J19").Value = 1241.4
$ws.Range("L19").Value = 1241.4
$ws.Range("N19").Value = -1591.4
$ws.Range("H38").Value = 115.86667
$ws.Range("J38").Value = 336
$ws.Range("L38").Value = 1008
$ws.Range("N38").Value = -1752
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H70").Value = 4537.5
$ws.Range("I70").Value = 3900
$ws.Range("J70").Value = 4628.5713
$ws.Range("K70").Value = 11700
$ws.Range("L70").Value = 13885.7139
$ws.Range("M70").Value = -11430
$ws.Range("N70").Value = -14425.7139
$ws.Range("H73").Value = 4537.5
$ws.Range("I73").Value = 3900
$ws.Range("J73").Value = 4628.5713
$ws.Range("K73").Value = 11700
$ws.Range("L73").Value = 13885.7139
$ws.Range("M73").Value = -10764
$ws.Range("N73").Value = -15757.7139
$ws.Range("H113").Value = 81570.57000000001
$ws.Range("I113").Value = 147998.72
$ws.Range("J113").Value = 15142.429
$ws.Range("K113").Value = 147998.72
$ws.Range("L113").Value = 15142.429
$ws.Range("M113").Value = -144744.72
$ws.Range("N113").Value = -21650.429
$ws.Range("H137").Value = 1874.2
$ws.Range("I137").Value = 1731
$ws.Range("K137").Value = 5193
$ws.Range("M137").Value = -2643
$ws.Range("H138").Value = 4074.0356
$ws.Range("J138").Value = 4582.5654
$ws.Range("L138").Value = 13747.6962
$ws.Range("N138").Value = -24027.6962

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3818.2307
$ws.Range("I61").Value = 3038.5
$ws.Range("K61").Value = 3038.5
$ws.Range("M61").Value = -2826.5
$ws.Range("H74").Value = 21913.777
$ws.Range("I74").Value = 2009.6666
$ws.Range("J74").Value = 31865.834
$ws.Range("K74").Value = 2009.6666
$ws.Range("L74").Value = 31865.834
$ws.Range("M74").Value = -1135.6666
$ws.Range("N74").Value = -33613.834
$ws.Range("H77").Value = 21913.777
$ws.Range("I77").Value = 2009.6666
$ws.Range("J77").Value = 31865.834
$ws.Range("K77").Value = 10048.333
$ws.Range("L77").Value = 159329.17
$ws.Range("M77").Value = -5680.333000000001
$ws.Range("N77").Value = -168065.17
$ws.Range("H132").Value = 3412.2942
$ws.Range("I132").Value = 3139.3
$ws.Range("J132").Value = 3802.2856
$ws.Range("K132").Value = 9417.900000000001
$ws.Range("L132").Value = 11406.8568
$ws.Range("M132").Value = -6887.900000000001
$ws.Range("N132").Value = -16466.8568
$ws.Range("H136").Value = 3818.2307
$ws.Range("I136").Value = 3038.5
$ws.Range("K136").Value = 9115.5
$ws.Range("M136").Value = -6565.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1419.4546
$ws.Range("I16").Value = 1213.625
$ws.Range("K16").Value = 1213.625
$ws.Range("M16").Value = -926.625
$ws.Range("H31").Value = 41184.703
$ws.Range("I31").Value = 64292.375
$ws.Range("K31").Value = 64292.375
$ws.Range("M31").Value = -63997.375
$ws.Range("H34").Value = 41184.703
$ws.Range("I34").Value = 64292.375
$ws.Range("K34").Value = 64292.375
$ws.Range("M34").Value = -64090.375
$ws.Range("H58").Value = 2184.158
$ws.Range("J58").Value = 1347
$ws.Range("L58").Value = 1347
$ws.Range("N58").Value = -1753
$ws.Range("H62").Value = 9687.875
$ws.Range("I62").Value = 8333.833000000001
$ws.Range("J62").Value = 13750
$ws.Range("K62").Value = 8333.833000000001
$ws.Range("L62").Value = 13750
$ws.Range("M62").Value = -7709.833000000001
$ws.Range("N62").Value = -14998
$ws.Range("H65").Value = 9687.875
$ws.Range("I65").Value = 8333.833000000001
$ws.Range("J65").Value = 13750
$ws.Range("K65").Value = 41669.165
$ws.Range("L65").Value = 68750
$ws.Range("M65").Value = -38549.165
$ws.Range("N65").Value = -74990
$ws.Range("H113").Value = 1419.4546
$ws.Range("I113").Value = 1213.625
$ws.Range("K113").Value = 1213.625
$ws.Range("M113").Value = 956.375
$ws.Range("H122").Value = 1788.8
$ws.Range("I122").Value = 1482.5
$ws.Range("K122").Value = 4447.5
$ws.Range("M122").Value = -1997.5
$ws.Range("H136").Value = 2184.158
$ws.Range("J136").Value = 1347
$ws.Range("L136").Value = 4041
$ws.Range("N136").Value = -9141

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 312.9524
$ws.Range("I40").Value = 199.27272
$ws.Range("K40").Value = 797.09088
$ws.Range("M40").Value = -728.09088

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 7999.8335
$ws.Range("J35").Value = 7999.8335
$ws.Range("L35").Value = 7999.8335
$ws.Range("N35").Value = -8595.833500000001
$ws.Range("H42").Value = 200412.14
$ws.Range("J42").Value = 200577
$ws.Range("L42").Value = 200577
$ws.Range("N42").Value = -201547
$ws.Range("H43").Value = 15104.25
$ws.Range("I43").Value = 15119.143
$ws.Range("K43").Value = 15119.143
$ws.Range("M43").Value = -14968.143
$ws.Range("H70").Value = 10593.6
$ws.Range("I70").Value = 8797.799999999999
$ws.Range("J70").Value = 11491.5
$ws.Range("K70").Value = 8797.799999999999
$ws.Range("L70").Value = 11491.5
$ws.Range("M70").Value = -8527.799999999999
$ws.Range("N70").Value = -12031.5
$ws.Range("H73").Value = 10593.6
$ws.Range("I73").Value = 8797.799999999999
$ws.Range("J73").Value = 11491.5
$ws.Range("K73").Value = 8797.799999999999
$ws.Range("L73").Value = 11491.5
$ws.Range("M73").Value = -7861.799999999999
$ws.Range("N73").Value = -13363.5
$ws.Range("H80").Value = 5203.4
$ws.Range("J80").Value = 5504.25
$ws.Range("L80").Value = 5504.25
$ws.Range("N80").Value = -7500.25
$ws.Range("H83").Value = 5203.4
$ws.Range("J83").Value = 5504.25
$ws.Range("L83").Value = 27521.25
$ws.Range("N83").Value = -37505.25
$ws.Range("H107").Value = 4000
$ws.Range("I107").Value = 4000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 4000
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -2080
$ws.Range("N107").ClearContents()
$ws.Range("H115").Value = 200412.14
$ws.Range("J115").Value = 200577
$ws.Range("L115").Value = 200577
$ws.Range("N115").Value = -202927

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 930.7
$ws.Range("I22").Value = 563
$ws.Range("J22").Value = 1175.8334
$ws.Range("K22").Value = 563
$ws.Range("L22").Value = 1175.8334
$ws.Range("M22").Value = -268
$ws.Range("N22").Value = -1765.8334
$ws.Range("H27").Value = 930.7
$ws.Range("I27").Value = 563
$ws.Range("J27").Value = 1175.8334
$ws.Range("K27").Value = 563
$ws.Range("L27").Value = 1175.8334
$ws.Range("M27").Value = -456
$ws.Range("N27").Value = -1389.8334
$ws.Range("H40").Value = 4734.5713
$ws.Range("I40").Value = 4023.25
$ws.Range("J40").Value = 9002.5
$ws.Range("K40").Value = 4023.25
$ws.Range("L40").Value = 9002.5
$ws.Range("M40").Value = -3887.25
$ws.Range("N40").Value = -9274.5
$ws.Range("H55").Value = 301.5
$ws.Range("I55").Value = 321.5
$ws.Range("J55").Value = 261.5
$ws.Range("K55").Value = 321.5
$ws.Range("L55").Value = 261.5
$ws.Range("M55").Value = -148.5
$ws.Range("N55").Value = -607.5
$ws.Range("H122").Value = 4921.4
$ws.Range("I122").Value = 4182.6
$ws.Range("K122").Value = 12547.8
$ws.Range("M122").Value = -10097.8
$ws.Range("H132").Value = 6510.5
$ws.Range("I132").Value = 5750
$ws.Range("J132").Value = 7017.5
$ws.Range("K132").Value = 17250
$ws.Range("L132").Value = 21052.5
$ws.Range("M132").Value = -14720
$ws.Range("N132").Value = -26112.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2423
$ws.Range("I122").Value = 2267.0513
$ws.Range("J122").Value = 2803.125
$ws.Range("K122").Value = 6801.1539
$ws.Range("L122").Value = 8409.375
$ws.Range("M122").Value = -4351.1539
$ws.Range("N122").Value = -13309.375
$ws.Range("H125").Value = 29999.2
$ws.Range("J125").Value = 29999.2
$ws.Range("L125").Value = 29999.2
$ws.Range("N125").Value = -39839.2
$ws.Range("H132").Value = 3369.2222
$ws.Range("I132").Value = 3172.3809
$ws.Range("K132").Value = 9517.1427
$ws.Range("M132").Value = -6987.1427
